# Replace the 25 two-digit-multiplication answers in the single 5-column
# table. The table has 20 rows; only rows 1, 5, 10, 15 and 20 contain the
# answer cells (the other rows are spacer rows). Cells are addressed
# directly by (row, column) rather than via text Find/Replace because one
# of the new values ("60x80=4800") is identical to an old value that is
# itself being replaced elsewhere in the table, which would make a plain
# global Find/Replace ambiguous / order-dependent.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("60×80=4800", "76×85=6460", "42×29=1218", "26×55=1430", "44×90=3960")
    5  = @("48×39=1872", "45×52=2340", "40×78=3120", "87×97=8439", "44×43=1892")
    10 = @("87×44=3828", "87×99=8613", "35×67=2345", "69×19=1311", "39×43=1677")
    15 = @("43×86=3698", "99×30=2970", "96×45=4320", "65×78=5070", "83×87=7221")
    20 = @("38×36=1368", "19×64=1216", "43×43=1849", "68×61=4148", "69×78=5382")
}

foreach ($rowIndex in $newValues.Keys) {
    $rowValues = $newValues[$rowIndex]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
